# Stop selecting the "WithUnicode" tab and instead finish with the new
# "ExtraBlankRowsAfterData" tab active - also append the two new sheets
# ("ExtraBlankRowsAfterData" and "Sheet2") at the end of the workbook so
# that no more "extra" blank rows linger after the real data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Sheet: ExtraBlankRowsAfterData -----------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "ExtraBlankRowsAfterData"

$ws3.Range("A1").Value = "aa11"
$ws3.Range("B1").Value = "bb22"
$ws3.Range("A2").Value = "cc33"
$ws3.Range("B2").Value = "dd44"

# A1 carries a custom numeric format.
$ws3.Range("A1").NumberFormat = "#.00,"

# Rows 3:4 (A:C) and row 5 (A:C) are leftover blank - but styled - rows
# that trail the real data (the "extra blank rows" this change targets).
$ws3.Range("A3:C4").Font.Color = 0
$ws3.Range("A3:C4").NumberFormat = "0.0?"
$ws3.Range("A5:C5").Font.Color = 0

$ws3.Range("A1:B2").Select()

# --- Sheet: Sheet2 (the clean, no-trailing-rows version) --------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Sheet2"

$ws4.Range("A1").Value = "aa11"
$ws4.Range("B1").Value = "bb22"
$ws4.Range("A2").Value = "cc33"
$ws4.Range("B2").Value = "dd44"

$ws4.Range("A1").NumberFormat = "#.00,"

$ws4.Range("A1:B2").Select()

# Leave "ExtraBlankRowsAfterData" as the active/visible tab.
$ws3.Activate()

# Shrink the saved window height slightly, matching the author's resize.
$excel.ActiveWindow.Height = 16400
